# Add 2022-Q4 data:
#  - insert a new "2022-Q4" sheet (cloned from the existing "2022-Q2" sheet so
#    it inherits identical formatting/styles) positioned right after "总计"
#    and before "2022-Q2"
#  - populate it with the 2022-Q4 fund holdings
#  - add a matching summary row on the "总计" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 with the 2022-Q4 summary figures
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.04

# A2 should carry the same style as the other index cells (A3:A6)
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0

# ---------------------------------------------------------------------------
# 2. Create the "2022-Q4" worksheet by cloning "2022-Q2" (same column layout
#    and styling), then overwrite its values with the 2022-Q4 fund data.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(2)
$templateSheet.Copy($null, $wb.Worksheets.Item(1))

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# the template only has one data row (row 2) - add 3 more for the 4 funds
$q4.Rows.Item(3).Insert()
$q4.Rows.Item(3).Insert()
$q4.Rows.Item(3).Insert()

# give the new index cells (A3:A5) the same style as A2
$q4.Range("A2").Copy()
$q4.Range("A3:A5").PasteSpecial(-4122)

# -- row 2: 007251 广发睿享稳健增利混合A --
$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "007251"
$q4.Range("C2").Value = "广发睿享稳健增利混合A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "1.26"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "39.64"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "1.73"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0218"
$q4.Range("H2").Value = 9

# -- row 3: 002137 诺安利鑫灵活配置混合A --
$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "002137"
$q4.Range("C3").Value = "诺安利鑫灵活配置混合A"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.44"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "89.87"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "3.42"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0150"
$q4.Range("H3").Value = 9

# -- row 4: 014521 诺安利鑫灵活配置混合C --
$q4.Range("A4").Value = 2
$q4.Range("B4").NumberFormat = "@"
$q4.Range("B4").Value = "014521"
$q4.Range("C4").Value = "诺安利鑫灵活配置混合C"
$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "0.01"
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "89.87"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "3.42"
$q4.Range("G4").NumberFormat = "@"
$q4.Range("G4").Value = "0.0003"
$q4.Range("H4").Value = 9

# -- row 5: 011702 广发睿享稳健增利混合C --
$q4.Range("A5").Value = 3
$q4.Range("B5").NumberFormat = "@"
$q4.Range("B5").Value = "011702"
$q4.Range("C5").Value = "广发睿享稳健增利混合C"
$q4.Range("D5").NumberFormat = "@"
$q4.Range("D5").Value = "0.00"
$q4.Range("E5").NumberFormat = "@"
$q4.Range("E5").Value = "39.64"
$q4.Range("F5").NumberFormat = "@"
$q4.Range("F5").Value = "1.73"
$q4.Range("G5").Value = 0
$q4.Range("H5").Value = 9
